$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.132183908045977
$ws.Range("C2").Value = 0.6417624521072797
$ws.Range("J2").Value = 0.01340996168582376
$ws.Range("P2").Value = 0.09386973180076628
$ws.Range("S2").Value = 0.1187739463601533
$ws.Range("B3").Value = 0.005649717514124294
$ws.Range("C3").Value = 0.04519774011299435
$ws.Range("J3").Value = 0.01977401129943503
$ws.Range("P3").Value = 0.7824858757062146
$ws.Range("S3").Value = 0.1468926553672316
$ws.Range("J4").Value = 0.05825242718446602
$ws.Range("P4").Value = 0.6213592233009708
$ws.Range("S4").Value = 0.3203883495145631
$ws.Range("B6").Value = 0.04826254826254826
$ws.Range("D6").Value = 0.01158301158301158
$ws.Range("F6").Value = 0.0694980694980695
$ws.Range("J6").Value = 0.2277992277992278
$ws.Range("O6").Value = 0.01158301158301158
$ws.Range("Q6").Value = 0.1795366795366795
$ws.Range("R6").Value = 0.06563706563706563
$ws.Range("S6").Value = 0.3861003861003861
$ws.Range("B7").Value = 0.1128668171557562
$ws.Range("D7").Value = 0.03386004514672687
$ws.Range("F7").Value = 0.04740406320541761
$ws.Range("J7").Value = 0.1399548532731377
$ws.Range("O7").Value = 0.02257336343115124
$ws.Range("Q7").Value = 0.1602708803611738
$ws.Range("R7").Value = 0.08803611738148984
$ws.Range("S7").Value = 0.3950338600451467
$ws.Range("B8").Value = 0.07516650808753568
$ws.Range("D8").Value = 0.02093244529019981
$ws.Range("E8").Value = 0.0009514747859181732
$ws.Range("F8").Value = 0.08372978116079924
$ws.Range("J8").Value = 0.08943862987630828
$ws.Range("O8").Value = 0.01712654614652712
$ws.Range("Q8").Value = 0.1845861084681256
$ws.Range("R8").Value = 0.1132254995242626
$ws.Range("S8").Value = 0.4148430066603235
$ws.Range("B9").Value = 0.09065934065934066
$ws.Range("D9").Value = 0.01373626373626374
$ws.Range("E9").Value = 0.002747252747252747
$ws.Range("F9").Value = 0.08241758241758242
$ws.Range("J9").Value = 0.07142857142857142
$ws.Range("O9").Value = 0.02197802197802198
$ws.Range("Q9").Value = 0.1758241758241758
$ws.Range("R9").Value = 0.1208791208791209
$ws.Range("S9").Value = 0.4203296703296703
$ws.Range("B10").Value = 0.1017871017871018
$ws.Range("D10").Value = 0.02292152292152292
$ws.Range("E10").Value = 0.001554001554001554
$ws.Range("F10").Value = 0.07964257964257965
$ws.Range("J10").Value = 0.09207459207459208
$ws.Range("O10").Value = 0.01515151515151515
$ws.Range("Q10").Value = 0.2257187257187257
$ws.Range("R10").Value = 0.09324009324009325
$ws.Range("S10").Value = 0.3679098679098679
$ws.Range("G11").Value = 0.1710914454277286
$ws.Range("J11").Value = 0.06932153392330384
$ws.Range("K11").Value = 0.1976401179941003
$ws.Range("L11").Value = 0.5486725663716814
$ws.Range("S11").Value = 0.01327433628318584
$ws.Range("G12").Value = 0.7413793103448276
$ws.Range("J12").Value = 0.1798029556650246
$ws.Range("K12").Value = 0.007389162561576354
$ws.Range("L12").Value = 0.03694581280788178
$ws.Range("S12").Value = 0.03448275862068965
$ws.Range("F13").Value = 0.01834862385321101
$ws.Range("G13").Value = 0.4954128440366973
$ws.Range("J13").Value = 0.4036697247706422
$ws.Range("S13").Value = 0.08256880733944955
$ws.Range("F15").Value = 0.02191235059760956
$ws.Range("H15").Value = 0.1852589641434263
$ws.Range("I15").Value = 0.06374501992031872
$ws.Range("J15").Value = 0.3725099601593626
$ws.Range("K15").Value = 0.06772908366533864
$ws.Range("M15").Value = 0.01195219123505976
$ws.Range("O15").Value = 0.05378486055776893
$ws.Range("S15").Value = 0.2231075697211155
$ws.Range("F16").Value = 0.01305483028720627
$ws.Range("H16").Value = 0.185378590078329
$ws.Range("I16").Value = 0.04960835509138381
$ws.Range("J16").Value = 0.4203655352480418
$ws.Range("K16").Value = 0.1201044386422977
$ws.Range("M16").Value = 0.02088772845953003
$ws.Range("N16").Value = 0.002610966057441253
$ws.Range("O16").Value = 0.06527415143603134
$ws.Range("S16").Value = 0.1227154046997389
$ws.Range("F17").Value = 0.01703406813627254
$ws.Range("H17").Value = 0.2074148296593187
$ws.Range("I17").Value = 0.0781563126252505
$ws.Range("J17").Value = 0.4228456913827655
$ws.Range("K17").Value = 0.08216432865731463
$ws.Range("M17").Value = 0.01803607214428858
$ws.Range("N17").Value = 0.001002004008016032
$ws.Range("O17").Value = 0.06513026052104208
$ws.Range("S17").Value = 0.1082164328657315
$ws.Range("F18").Value = 0.008456659619450317
$ws.Range("H18").Value = 0.2114164904862579
$ws.Range("I18").Value = 0.07610993657505286
$ws.Range("J18").Value = 0.3932346723044398
$ws.Range("K18").Value = 0.1162790697674419
$ws.Range("M18").Value = 0.03594080338266385
$ws.Range("O18").Value = 0.06765327695560254
$ws.Range("S18").Value = 0.09090909090909091
$ws.Range("F19").Value = 0.02066590126291619
$ws.Range("H19").Value = 0.2166092613853808
$ws.Range("I19").Value = 0.07615767317259854
$ws.Range("J19").Value = 0.3616532721010333
$ws.Range("K19").Value = 0.1224646000765404
$ws.Range("M19").Value = 0.02257941063911213
$ws.Range("N19").Value = 0.0007654037504783774
$ws.Range("O19").Value = 0.07921928817451206
$ws.Range("S19").Value = 0.09988518943742825
